$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.705.88"
$ws.Range("E2").Value = "  +3.74%  "

$ws.Range("D3").Value = "2.417.29"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.87%  "

$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.68%  "

$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("E12").Value = "  -2.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").Value = "2.794.66"
$ws.Range("E15").Value = "  +2.06%  "

$ws.Range("D16").Value = "2.399.93"
$ws.Range("E16").Value = "  +3.61%  "

$ws.Range("E17").Value = "  +3.54%  "

$ws.Range("D18").Value = "44.527.26"
$ws.Range("E18").Value = "  +3.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").Value = "0.0₃0914"
$ws.Range("E21").Value = "  +3.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.83%  "

$ws.Range("E24").Value = "  +3.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.16%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.65%  "

$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  +15.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0762"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "127.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.67%  "

$ws.Range("E42").Value = "  -3.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0287"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.99%  "

$ws.Range("D45").Value = "1.933.67"
$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("E47").Value = "  +6.94%  "

$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  +16.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.89%  "
